# Revisão para a avaliação
# Slide 10: split the opening run of the content placeholder paragraph into
# three runs: "Para " / "declarar " / "uma variável, primeiro é preciso
# definir qual o " (same run formatting throughout).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$para1 = $tr.Paragraphs(1, 1)
$run1 = $para1.Runs(1, 1)

# Shrink the existing run down to its trailing text…
$run1.Text = "uma variável, primeiro é preciso definir qual o "

# …then insert the two new leading runs in front of it (InsertBefore carries
# over the same run-level formatting from $run1).
$null = $run1.InsertBefore("declarar ")
$null = $run1.InsertBefore("Para ")
